$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2714.125
$ws.Range("I9").Value = 553.25
$ws.Range("J9").Value = 4875
$ws.Range("K9").Value = 553.25
$ws.Range("L9").Value = 4875
$ws.Range("M9").Value = -384.25
$ws.Range("N9").Value = -5213
$ws.Range("H15").Value = 605.0877
$ws.Range("I15").Value = 605.0877
$ws.Range("K15").Value = 1815.2631
$ws.Range("M15").Value = -1646.2631
$ws.Range("H28").Value = 1013.38464
$ws.Range("I28").Value = 960.7273
$ws.Range("K28").Value = 960.7273
$ws.Range("M28").Value = -475.7273
$ws.Range("H64").Value = 11221.556
$ws.Range("I64").Value = 7165.8335
$ws.Range("K64").Value = 7165.8335
$ws.Range("M64").Value = -6917.8335
$ws.Range("H67").Value = 11221.556
$ws.Range("I67").Value = 7165.8335
$ws.Range("K67").Value = 7165.8335
$ws.Range("M67").Value = -6307.8335
$ws.Range("H70").Value = 4000
$ws.Range("J70").Value = 4000
$ws.Range("L70").Value = 12000
$ws.Range("N70").Value = -12540
$ws.Range("H73").Value = 4000
$ws.Range("J73").Value = 4000
$ws.Range("L73").Value = 12000
$ws.Range("N73").Value = -13872
$ws.Range("H107").Value = 766.0909
$ws.Range("I107").Value = 603.7143
$ws.Range("J107").Value = 1050.25
$ws.Range("K107").Value = 603.7143
$ws.Range("L107").Value = 1050.25
$ws.Range("M107").Value = 1316.2857
$ws.Range("N107").Value = -4890.25
$ws.Range("H137").Value = 3050
$ws.Range("I137").Value = 775
$ws.Range("J137").Value = 4187.5
$ws.Range("K137").Value = 2325
$ws.Range("L137").Value = 12562.5
$ws.Range("M137").Value = 225
$ws.Range("N137").Value = -17662.5
$ws.Range("H138").Value = 9002.223
$ws.Range("J138").Value = 9742.186
$ws.Range("L138").Value = 29226.558
$ws.Range("N138").Value = -39506.558

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 878
$ws.Range("I132").Value = 878
$ws.Range("K132").Value = 2634
$ws.Range("M132").Value = -104

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1606.3158
$ws.Range("J94").Value = 5666.3335
$ws.Range("L94").Value = 5666.3335
$ws.Range("N94").Value = -6568.3335
$ws.Range("H107").Value = 1010.087
$ws.Range("I107").Value = 1038.5264
$ws.Range("J107").Value = 875
$ws.Range("K107").Value = 1038.5264
$ws.Range("L107").Value = 875
$ws.Range("M107").Value = 881.4736
$ws.Range("N107").Value = -4715
$ws.Range("H134").Value = 2497.8462
$ws.Range("I134").Value = 2322.7
$ws.Range("J134").Value = 3081.6667
$ws.Range("K134").Value = 6968.099999999999
$ws.Range("L134").Value = 9245.000100000001
$ws.Range("M134").Value = -4433.099999999999
$ws.Range("N134").Value = -14315.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3149.7273
$ws.Range("I31").Value = 3102
$ws.Range("J31").Value = 3277
$ws.Range("K31").Value = 3102
$ws.Range("L31").Value = 3277
$ws.Range("M31").Value = -2807
$ws.Range("N31").Value = -3867
$ws.Range("H34").Value = 3149.7273
$ws.Range("I34").Value = 3102
$ws.Range("J34").Value = 3277
$ws.Range("K34").Value = 3102
$ws.Range("L34").Value = 3277
$ws.Range("M34").Value = -2900
$ws.Range("N34").Value = -3681
$ws.Range("H99").Value = 7320.846
$ws.Range("I99").Value = 7128.6665
$ws.Range("J99").Value = 7753.25
$ws.Range("K99").Value = 7128.6665
$ws.Range("L99").Value = 7753.25
$ws.Range("M99").Value = -5630.6665
$ws.Range("N99").Value = -10749.25
$ws.Range("H126").Value = 7320.846
$ws.Range("I126").Value = 7128.6665
$ws.Range("J126").Value = 7753.25
$ws.Range("K126").Value = 21385.9995
$ws.Range("L126").Value = 23259.75
$ws.Range("M126").Value = -18915.9995
$ws.Range("N126").Value = -28199.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 280.55554
$ws.Range("I50").Value = 282.35294
$ws.Range("J50").Value = 250
$ws.Range("K50").Value = 847.05882
$ws.Range("L50").Value = 750
$ws.Range("M50").Value = -366.05882
$ws.Range("N50").Value = -1712
$ws.Range("H53").Value = 280.55554
$ws.Range("I53").Value = 282.35294
$ws.Range("J53").Value = 250
$ws.Range("K53").Value = 847.05882
$ws.Range("L53").Value = 750
$ws.Range("M53").Value = -366.05882
$ws.Range("N53").Value = -1712
$ws.Range("H140").Value = 1574.75
$ws.Range("I140").Value = 599.6667
$ws.Range("K140").Value = 1799.0001
$ws.Range("M140").Value = 3380.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()
$ws.Range("H80").Value = 3364.3333
$ws.Range("I80").Value = 3300
$ws.Range("K80").Value = 3300
$ws.Range("M80").Value = -2302
$ws.Range("H83").Value = 3364.3333
$ws.Range("I83").Value = 3300
$ws.Range("K83").Value = 16500
$ws.Range("M83").Value = -11508
$ws.Range("H122").Value = 3355.2273
$ws.Range("I122").Value = 3361.1765
$ws.Range("J122").Value = 3335
$ws.Range("K122").Value = 10083.5295
$ws.Range("L122").Value = 10005
$ws.Range("M122").Value = -7633.529500000001
$ws.Range("N122").Value = -14905
$ws.Range("H127").Value = 49999
$ws.Range("J127").Value = 49999
$ws.Range("L127").Value = 49999
$ws.Range("N127").Value = -59919

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 23329.666
$ws.Range("J40").Value = 20000
$ws.Range("L40").Value = 20000
$ws.Range("N40").Value = -20272
$ws.Range("H82").Value = 3384.9
$ws.Range("I82").Value = 3098.6667
$ws.Range("J82").Value = 3507.5715
$ws.Range("K82").Value = 3098.6667
$ws.Range("L82").Value = 3507.5715
$ws.Range("M82").Value = -2737.6667
$ws.Range("N82").Value = -4229.5715
$ws.Range("H85").Value = 3384.9
$ws.Range("I85").Value = 3098.6667
$ws.Range("J85").Value = 3507.5715
$ws.Range("K85").Value = 3098.6667
$ws.Range("L85").Value = 3507.5715
$ws.Range("M85").Value = -1850.6667
$ws.Range("N85").Value = -6003.5715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1259.6
$ws.Range("H84").Value = 1259.6
$ws.Range("H126").Value = 1654.9286
$ws.Range("I126").Value = 1308.75
$ws.Range("K126").Value = 3926.25
$ws.Range("M126").Value = -1456.25
